$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "LP1912"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 02:55:01"
$ws1.Range("A3").Value = "Total filas: 12"

# Row 10 updated
$ws1.Range("A10").Value = "02:55:01"
$ws1.Range("B10").Value = "02:58"
$ws1.Range("C10").Value = "215_ALUAR"
$ws1.Range("D10").Value = 3
$ws1.Range("E10").Value = "LP1912"

# Row 11 stays the same (01:59:40 | 03:06 | 215_ALUAR | 67 | LP1912)

# Row 12 updated
$ws1.Range("A12").Value = "02:55:01"
$ws1.Range("B12").Value = "03:48"
$ws1.Range("C12").Value = "14_ABASTO"
$ws1.Range("D12").Value = 53
$ws1.Range("E12").Value = "LP1912"

# Row 13 updated
$ws1.Range("A13").Value = "01:59:40"
$ws1.Range("B13").Value = "03:50"
$ws1.Range("C13").Value = "14_ABASTO"
$ws1.Range("D13").Value = 111
$ws1.Range("E13").Value = "LP1912"

# Row 14 updated
$ws1.Range("A14").Value = "02:30:53"
$ws1.Range("B14").Value = "03:52"
$ws1.Range("C14").Value = "14_ABASTO"
$ws1.Range("D14").Value = 82
$ws1.Range("E14").Value = "LP1912"

# Row 15 (new)
$ws1.Range("A15").Value = "02:55:01"
$ws1.Range("B15").Value = "04:01"
$ws1.Range("C15").Value = "81_EL PELIGRO"
$ws1.Range("D15").Value = 66
$ws1.Range("E15").Value = "LP1912"

# Row 16 (new)
$ws1.Range("A16").Value = "02:55:01"
$ws1.Range("B16").Value = "04:46"
$ws1.Range("C16").Value = "215A_EL PATO"
$ws1.Range("D16").Value = 111
$ws1.Range("E16").Value = "LP1912"

# Row 17 (new)
$ws1.Range("A17").Value = "02:55:01"
$ws1.Range("B17").Value = "04:53"
$ws1.Range("C17").Value = "11_ETCHEVERRY"
$ws1.Range("D17").Value = 118
$ws1.Range("E17").Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet "LP1912-215"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 02:55:01"
$ws2.Range("A3").Value = "Total filas: 5"

# Row 8 updated
$ws2.Range("A8").Value = "02:55:01"
$ws2.Range("B8").Value = "02:58"
$ws2.Range("C8").Value = "215_ALUAR"
$ws2.Range("D8").Value = 3
$ws2.Range("E8").Value = "LP1912"

# Row 9 stays the same (01:59:40 | 03:06 | 215_ALUAR | 67 | LP1912)

# Row 10 (new)
$ws2.Range("A10").Value = "02:55:01"
$ws2.Range("B10").Value = "04:46"
$ws2.Range("C10").Value = "215A_EL PATO"
$ws2.Range("D10").Value = 111
$ws2.Range("E10").Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet "6203-6173"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 02:55:01"
